$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-113 down to 31-114.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44414
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112008
$ws.Range("G30").Value = "Coliflor"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Segunda"
$ws.Range("J30").Value = 3000
$ws.Range("K30").Value = 500
$ws.Range("L30").Value = 500
$ws.Range("M30").Value = 500
$ws.Range("N30").Value = '$/unidad'
$ws.Range("O30").Value = "Región del Maule"
$ws.Range("P30").Value = 500
$ws.Range("Q30").Value = 1
$ws.Range("R30").Value = "Hortaliza"
